$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Light-blue fill (the workbook's existing "in progress" style, fill color FF00B0F0)
# expressed as a BGR-packed OLE color value for Interior.Color.
$blueColor = 15773696

# Contiguous A:C row blocks whose instruction rows are now marked with the
# blue "in progress" background instead of the default (no-fill) style.
$rowRanges = @(
    @(12, 15),
    @(17, 17),
    @(20, 20),
    @(22, 46),
    @(70, 88),
    @(93, 103)
)

foreach ($rr in $rowRanges) {
    $startRow = $rr[0]
    $endRow = $rr[1]
    $range = $ws.Range("A$startRow`:C$endRow")
    $range.Interior.Color = $blueColor
}

# Move the active selection/view down to the bottom of the table (A93:C103)
$ws.Range("A93:C103").Select()
